# Applies the scheduled market-data refresh to Sheets H:N price/profit columns.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 224.35
$ws.Range("I39").Value = 246.41176
$ws.Range("J39").Value = 99.333336
$ws.Range("K39").Value = 739.23528
$ws.Range("L39").Value = 298.000008
$ws.Range("M39").Value = -443.23528
$ws.Range("N39").Value = -890.000008
# Row 62
$ws.Range("H62").Value = 29509.531
$ws.Range("I62").Value = 24844.215
$ws.Range("J62").Value = 33138.11
$ws.Range("K62").Value = 24844.215
$ws.Range("L62").Value = 33138.11
$ws.Range("M62").Value = -24220.215
$ws.Range("N62").Value = -34386.11
# Row 65
$ws.Range("H65").Value = 29509.531
$ws.Range("I65").Value = 24844.215
$ws.Range("J65").Value = 33138.11
$ws.Range("K65").Value = 124221.075
$ws.Range("L65").Value = 165690.55
$ws.Range("M65").Value = -121101.075
$ws.Range("N65").Value = -171930.55
# Row 80
$ws.Range("H80").Value = 1137498.8
$ws.Range("I80").Value = 2273477
$ws.Range("J80").Value = 1520.5
$ws.Range("K80").Value = 6820431
$ws.Range("L80").Value = 4561.5
$ws.Range("M80").Value = -6819433
$ws.Range("N80").Value = -6557.5
# Row 83
$ws.Range("H83").Value = 1137498.8
$ws.Range("I83").Value = 2273477
$ws.Range("J83").Value = 1520.5
$ws.Range("K83").Value = 20461293
$ws.Range("L83").Value = 13684.5
$ws.Range("M83").Value = -20456301
$ws.Range("N83").Value = -23668.5
# Row 86
$ws.Range("H86").Value = 11143919
$ws.Range("I86").Value = 5819.875
$ws.Range("J86").Value = 20054398
$ws.Range("K86").Value = 5819.875
$ws.Range("L86").Value = 20054398
$ws.Range("M86").Value = -4696.875
$ws.Range("N86").Value = -20056644
# Row 89
$ws.Range("H89").Value = 11143919
$ws.Range("I89").Value = 5819.875
$ws.Range("J89").Value = 20054398
$ws.Range("K89").Value = 29099.375
$ws.Range("L89").Value = 100271990
$ws.Range("M89").Value = -23483.375
$ws.Range("N89").Value = -100283222
# Row 96
$ws.Range("H96").Value = 475.0909
$ws.Range("I96").Value = 303.7143
$ws.Range("K96").Value = 911.1428999999999
$ws.Range("M96").Value = 461.8571000000001
# Row 100
$ws.Range("H100").Value = 4407.294
$ws.Range("I100").Value = 4495.25
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 4495.25
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -3954.25
$ws.Range("N100").Value = -4082
# Row 112
$ws.Range("H112").Value = 47627.863
$ws.Range("J112").Value = 2095.4211
$ws.Range("L112").Value = 6286.263300000001
$ws.Range("N112").Value = -8502.263300000001
# Row 125
$ws.Range("H125").Value = 2250
$ws.Range("I125").Value = 1850
$ws.Range("K125").Value = 16650
$ws.Range("M125").Value = -14190
# Row 138
$ws.Range("H138").Value = 2965.718
$ws.Range("I138").Value = 2247.375
$ws.Range("J138").Value = 6249.5713
$ws.Range("K138").Value = 6742.125
$ws.Range("L138").Value = 18748.7139
$ws.Range("M138").Value = -1602.125
$ws.Range("N138").Value = -29028.7139

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 10048.315
$ws.Range("I102").Value = 7863
$ws.Range("K102").Value = 7863
$ws.Range("M102").Value = -6241
# Row 132
$ws.Range("H132").Value = 669659.7
$ws.Range("I132").Value = 457471.03
$ws.Range("K132").Value = 1372413.09
$ws.Range("M132").Value = -1369883.09

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 33335296
$ws.Range("I86").Value = 50001444
$ws.Range("J86").Value = 3002
$ws.Range("K86").Value = 50001444
$ws.Range("L86").Value = 3002
$ws.Range("M86").Value = -50000321
$ws.Range("N86").Value = -5248
# Row 89
$ws.Range("H89").Value = 33335296
$ws.Range("I89").Value = 50001444
$ws.Range("J89").Value = 3002
$ws.Range("K89").Value = 250007220
$ws.Range("L89").Value = 15010
$ws.Range("M89").Value = -250001604
$ws.Range("N89").Value = -26242
# Row 99
$ws.Range("H99").Value = 1458.6316
$ws.Range("I99").Value = 1458.6316
$ws.Range("K99").Value = 1458.6316
$ws.Range("M99").Value = 39.36840000000007

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4996.25
$ws.Range("I16").Value = 4995
$ws.Range("K16").Value = 4995
$ws.Range("M16").Value = -4708
# Row 22
$ws.Range("H22").Value = 1470.579
$ws.Range("I22").Value = 449.81818
$ws.Range("K22").Value = 449.81818
$ws.Range("M22").Value = -99.81817999999998
# Row 108
$ws.Range("H108").Value = 353334
$ws.Range("J108").Value = 353334
$ws.Range("L108").Value = 353334
$ws.Range("N108").Value = -361014
# Row 113
$ws.Range("H113").Value = 4996.25
$ws.Range("I113").Value = 4995
$ws.Range("K113").Value = 4995
$ws.Range("M113").Value = -2825
# Row 130
$ws.Range("H130").Value = 39406.668
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 167.92857
$ws.Range("I33").Value = 146.55556
$ws.Range("J33").Value = 206.4
$ws.Range("K33").Value = 879.3333600000001
$ws.Range("L33").Value = 1238.4
$ws.Range("M33").Value = -596.3333600000001
$ws.Range("N33").Value = -1804.4
# Row 68
$ws.Range("H68").Value = 1860.2858
$ws.Range("J68").Value = 2199.4443
$ws.Range("L68").Value = 6598.3329
$ws.Range("N68").Value = -8220.332900000001
# Row 71
$ws.Range("H71").Value = 1860.2858
$ws.Range("J71").Value = 2199.4443
$ws.Range("L71").Value = 19794.9987
$ws.Range("N71").Value = -27906.9987
# Row 131
$ws.Range("H131").Value = 4782.6665
$ws.Range("I131").Value = 914.5
$ws.Range("J131").Value = 8948.385
$ws.Range("K131").Value = 2743.5
$ws.Range("L131").Value = 26845.155
$ws.Range("M131").Value = 2296.5
$ws.Range("N131").Value = -36925.155
# Row 139
$ws.Range("H139").Value = 1456.625
$ws.Range("I139").Value = 1456.625
$ws.Range("K139").Value = 4369.875
$ws.Range("M139").Value = 770.125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 49999
$ws.Range("J26").Value = 49999
$ws.Range("L26").Value = 49999
$ws.Range("N26").Value = -50559
# Row 50
$ws.Range("H50").Value = 49999
$ws.Range("J50").Value = 49999
$ws.Range("L50").Value = 49999
$ws.Range("N50").Value = -50995
# Row 97
$ws.Range("H97").Value = 1315.4857
$ws.Range("I97").Value = 1315.2413
$ws.Range("K97").Value = 1315.2413
$ws.Range("M97").Value = -819.2412999999999
# Row 102
$ws.Range("H102").Value = 2393.7942
$ws.Range("I102").Value = 1741.8823
$ws.Range("J102").Value = 3045.7058
$ws.Range("K102").Value = 1741.8823
$ws.Range("L102").Value = 3045.7058
$ws.Range("M102").Value = -119.8823
$ws.Range("N102").Value = -6289.7058
# Row 132
$ws.Range("H132").Value = 230368.81
$ws.Range("I132").Value = 305630.38
$ws.Range("J132").Value = 4584.1816
$ws.Range("K132").Value = 916891.14
$ws.Range("L132").Value = 13752.5448
$ws.Range("M132").Value = -914361.14
$ws.Range("N132").Value = -18812.5448

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 27781166
$ws.Range("I7").Value = 62502564
$ws.Range("J7").Value = 4049.9
$ws.Range("K7").Value = 62502564
$ws.Range("L7").Value = 4049.9
$ws.Range("M7").Value = -62502452
$ws.Range("N7").Value = -4273.9
# Row 16
$ws.Range("H16").Value = 2461.2693
$ws.Range("I16").Value = 2074.375
$ws.Range("K16").Value = 2074.375
$ws.Range("M16").Value = -1904.375
# Row 40
$ws.Range("H40").Value = 3294.6
$ws.Range("I40").Value = 2834.7368
$ws.Range("J40").Value = 4750.8335
$ws.Range("K40").Value = 2834.7368
$ws.Range("L40").Value = 4750.8335
$ws.Range("M40").Value = -2698.7368
$ws.Range("N40").Value = -5022.8335
# Row 93
$ws.Range("H93").Value = 1505.4546
$ws.Range("I93").Value = 1603.75
$ws.Range("K93").Value = 1603.75
$ws.Range("M93").Value = -355.75
# Row 122
$ws.Range("H122").Value = 3848.1724
$ws.Range("I122").Value = 2914.0715
$ws.Range("J122").Value = 4720
$ws.Range("K122").Value = 8742.2145
$ws.Range("L122").Value = 14160
$ws.Range("M122").Value = -6292.2145
$ws.Range("N122").Value = -19060
# Row 126
$ws.Range("H126").Value = 27781166
$ws.Range("I126").Value = 62502564
$ws.Range("J126").Value = 4049.9
$ws.Range("K126").Value = 187507692
$ws.Range("L126").Value = 12149.7
$ws.Range("M126").Value = -187505222
$ws.Range("N126").Value = -17089.7
# Row 136
$ws.Range("H136").Value = 6926.933
$ws.Range("I136").Value = 2819.5
$ws.Range("K136").Value = 8458.5
$ws.Range("M136").Value = -5908.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Range("H61").Value = 3051
$ws.Range("I61").Value = 3051
$ws.Range("K61").Value = 3051
$ws.Range("M61").Value = -2759
# Row 109
$ws.Range("H109").Value = 79749.5
$ws.Range("J109").Value = 79749.5
$ws.Range("L109").Value = 79749.5
$ws.Range("N109").Value = -82523.5
# Row 122
$ws.Range("H122").Value = 1833.3
$ws.Range("I122").Value = 1391.875
$ws.Range("K122").Value = 4175.625
$ws.Range("M122").Value = -1725.625
# Row 126
$ws.Range("H126").Value = 1747.7742
$ws.Range("I126").Value = 1686.0435
$ws.Range("J126").Value = 1925.25
$ws.Range("K126").Value = 5058.1305
$ws.Range("L126").Value = 5775.75
$ws.Range("M126").Value = -2588.1305
$ws.Range("N126").Value = -10715.75
